$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.001983710570619874
$wsLinear.Range("B3").Value = -0.02970025362877381
$wsLinear.Range("B4").Value = 1.654071473114823
$wsLinear.Range("B5").Value = "[1.0, 0.22345764208746138, 0.009313186839776642, -0.04567286377204512, -0.04794271230765947, 0.02590096487705561, 0.19910422170768854, 0.3785768914540917, 0.19575101645769982, 0.013748199888822747, -0.04021901070608031, -0.06119875618416416, 0.014896822929713647, 0.19130792096365895, 0.3562147833776934, 0.17901019314583738, 0.0009082964563972481, -0.04525648710207334, -0.04514570313848705, -0.001331338540822674]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.0007870276248536576
$wsNonLinear.Range("B5").Value = -0.03749738162972732
$wsNonLinear.Range("B6").Value = 1.826142250295411
$wsNonLinear.Range("B7").Value = -0.009574743502662506
$wsNonLinear.Range("B8").Value = -0.006618982860083729
$wsNonLinear.Range("B9").Value = 1.500112687447744
$wsNonLinear.Range("B10").Value = "[0.9999999999999998, 0.22347918644848128, 0.009872236124642721, -0.0451442403432018, -0.0475061906717393, 0.026160453395892725, 0.1989688389438082, 0.3781145900016963, 0.195800814794052, 0.014181629877157362, -0.039720859714758076, -0.06075085612026895, 0.01527893871738951, 0.19126015588292286, 0.35590683097819037, 0.1789797633652115, 0.0012336470707487335, -0.044794967525722686, -0.044641799568667505, -0.0009625485340753756]"
